# "Subo ETL con datos actualizados"
# Append 10 newly-extracted records (ids 6-15) to the Sheet1 data table,
# which previously only held ids 1-5 (rows 2-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids      = @(6, 7, 8, 9, 10, 11, 12, 13, 14, 15)
$nombres  = @("María López", "Jorge Medina", "Rosa Castillo", "Luis Fernández", "Patricia Ríos", "Miguel Quispe", "Sandra Núñez", "Renzo Valdivia", "Carmen Soto", "Diego Paredes")
$edades   = @(31, 27, 42, 36, 29, 48, 33, 26, 41, 35)
$ciudades = @("Chiclayo", "Iquitos", "Huaraz", "Tacna", "Puno", "Ayacucho", "Huancayo", "Lima", "Moquegua", "Tumbes")
$importes = @(175, 130, 220, 195, 160, 240, 185, 140, 205, 190)

$firstNewRow = 7
$count = $ids.Length

# Write the new "nombre" values first, then the new "ciudad" values, so
# brand-new person names are registered ahead of the (mostly brand-new)
# city names, mirroring how the source ETL export ordered them.
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstNewRow + $i, 2).Value = $nombres[$i]
}
for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($firstNewRow + $i, 4).Value = $ciudades[$i]
}

# Fill in the remaining (numeric) columns for the new rows: id, edad, importe.
for ($i = 0; $i -lt $count; $i++) {
    $row = $firstNewRow + $i
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 3).Value = $edades[$i]
    $ws.Cells.Item($row, 5).Value = $importes[$i]
}
